# ---------------------------------------------------------------------------
# "Added technicians to printed data, fixed issue in case rho==0"
#
# Rewrites the "Year" worksheet (sheet1) of the RowCol-calc workbook:
#   - updates the four input parameters in B1:B4
#   - expands the "Rows" block (rows 12-19) with a new Deduction column (F)
#     and two new resource lines (Resources R Up / Resources R Low), plus
#     renames some rows, to avoid div-by-zero style issues (rho==0)
#   - expands the "Columns" block (rows 22-27) with two new technician
#     lines Mp / Mr
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Year")

# --- Clear cells that exist in the old layout but have no counterpart in
#     the new one (content AND formatting), so their shared strings / old
#     bold styling don't linger. ---------------------------------------
$ws.Range("G15").Clear()
$ws.Range("G16").Clear()
$ws.Range("G17").Clear()
$ws.Range("A20").Clear()
$ws.Range("B20").Clear()

# A18 used to be the bold "Columns" header (now A21) - drop the old bold
# formatting before it's reused for the plain "Maintenance 1" label.
$ws.Range("A18").Clear()

# G19 used to be a blank bold-styled cell - drop that formatting before
# it's reused for the plain "approx based on M and gL" label.
$ws.Range("G19").Clear()

# --- New shared strings must be introduced in the exact order they occupy
#     in the rebuilt sharedStrings table (Mp, Mr, Resources P Up/Low,
#     Resources R Low/Up, Deduction, approx based on M and gU/gL). ---------
$ws.Range("A26").Value = "Mp"
$ws.Range("A27").Value = "Mr"
$ws.Range("A13").Value = "Resources P Up"
$ws.Range("A14").Value = "Resources P Low"
$ws.Range("A16").Value = "Resources R Low"
$ws.Range("A15").Value = "Resources R Up"
$ws.Range("F11").Value = "Deduction"
$ws.Range("G18").Value = "approx based on M and gU"
$ws.Range("G19").Value = "approx based on M and gL"

# --- Input parameters (rows 1-4) ------------------------------------------
$ws.Range("B1").Value = 4
$ws.Range("B2").Value = 120
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 12

# --- Row 11: section headers / totals --------------------------------------
$ws.Range("B11").Formula = "=SUM(B12:B19)"
$ws.Range("E11").Formula = "=SUM(E12:E19)"
$ws.Range("F11").Font.Bold = $true

# --- Row 12 ------------------------------------------------------------
$ws.Range("A12").Value = "Capacity"
$ws.Range("B12").Formula = "=B3*B2*B1"
$ws.Range("D12").Value = 9.5
$ws.Range("E12").Formula = "=B12*D12-F12"
$ws.Range("F12").Value = 0

# --- Row 13: Resources P Up -------------------------------------------
$ws.Range("B13").Formula = "=B2*B1"
$ws.Range("D13").Value = 2
$ws.Range("E13").Formula = "=B13*D13-F13"
$ws.Range("F13").Value = 0

# --- Row 14: Resources P Low ------------------------------------------
$ws.Range("B14").Formula = "=B2*B1"
$ws.Range("D14").Value = 7
$ws.Range("E14").Formula = "=B14*D14-F14"
$ws.Range("F14").Value = 0

# --- Row 15: Resources R Up -------------------------------------------
$ws.Range("B15").Formula = "=B1*B2*B5"
$ws.Range("D15").Value = 2
$ws.Range("E15").Formula = "=B15*D15-F15"
$ws.Range("F15").Value = 0

# --- Row 16: Resources R Low ------------------------------------------
$ws.Range("B16").Formula = "=B1*B2*B3*B5"
$ws.Range("D16").Value = 1.5
$ws.Range("E16").Formula = "=B16*D16-F16"
$ws.Range("F16").Value = 0

# --- Row 17: Repair ------------------------------------------------------
$ws.Range("A17").Value = "Repair"
$ws.Range("B17").Formula = "=B3*B2*B5"
$ws.Range("D17").Value = 3
$ws.Range("E17").Formula = "=B17*D17-F17"
$ws.Range("F17").Formula = "=B3*B5"

# --- Row 18: Maintenance 1 ----------------------------------------------
$ws.Range("A18").Value = "Maintenance 1"
$ws.Range("B18").Formula = "=B4*B2"
$ws.Range("D18").Value = 105.5
$ws.Range("E18").Formula = "=B18*D18-F18"
$ws.Range("F18").Value = 0

# --- Row 19: Maintenance 2 ----------------------------------------------
$ws.Range("A19").Value = "Maintenance 2"
$ws.Range("B19").Formula = "=B4*B2"
$ws.Range("D19").Formula = "=D18"
$ws.Range("E19").Formula = "=B19*D19-F19"
$ws.Range("F19").Value = 0

# --- Row 21: "Columns" section header / total ---------------------------
$ws.Range("A21").Value = "Columns"
$ws.Range("A21").Font.Bold = $true
$ws.Range("B21").Formula = "=SUM(B22:B27)"

# --- Rows 22-27: column breakdown, incl. new Mp / Mr technician rows -----
$ws.Range("A22").Value = "N"
$ws.Range("B22").Formula = "=B1*B2"

$ws.Range("A23").Value = "P"
$ws.Range("B23").Formula = "=B2*B4"

$ws.Range("A24").Value = "R"
$ws.Range("B24").Formula = "=B2*B5*B3"

$ws.Range("A25").Value = "U"
$ws.Range("B25").Formula = "=B2*B5*B3"

$ws.Range("B26").Formula = "=B1*B2"

$ws.Range("B27").Formula = "=B1*B2*B5"

# --- Cosmetic layout tweaks mirrored from the authored commit -------------
$ws.Range("A1").Select()
$ws.Range("E11").Select()

$ws.Columns.Item(1).ColumnWidth = 15.7109375
$ws.Columns.Item(6).ColumnWidth = 10.140625

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
